# TI complet + TVML
# Adds the J5:S6 data series and a line chart (plotting rows 5 and 6)
# anchored below the existing data, matching the authored workbook edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: first data series (values used as the chart's first line) ---
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 4
$ws.Range("N5").Value = 4
$ws.Range("O5").Value = 7
$ws.Range("P5").Value = 7
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = 10
$ws.Range("S5").Value = 11

# --- Row 6: second data series ---
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0.2
$ws.Range("M6").Value = 0.2
$ws.Range("N6").Value = 0.1
$ws.Range("O6").Value = 0.1
$ws.Range("P6").Value = 0.4
$ws.Range("Q6").Value = 0.4
$ws.Range("R6").Value = 0.3
$ws.Range("S6").Value = 0.3

# --- Line chart fed from the two new rows, placed under the table (~A15:I32) ---
$chartObj = $ws.ChartObjects().Add(42, 189, 443.5, 216)
$chart = $chartObj.Chart
$chart.ChartType = 4   # xlLine

$chart.SeriesCollection().NewSeries()
$ser1 = $chart.SeriesCollection(1)
$ser1.Values = $ws.Range("J5:S5")

$ser2 = $chart.SeriesCollection(2)
$ser2.Values = $ws.Range("J6:S6")

# --- Restore the selection recorded in the authored edit ---
$ws.Range("R5:S6").Select()
